# Update the "varying alpha & m" sheet (active sheet) with refreshed
# simulation output for financially constrained firms (rows 34-51, cols N-X).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34
$ws.Range("N34").Value = [double]"7.3890961710000001E-3"
$ws.Range("O34").Value = [double]"6.9617070579999999E-3"
$ws.Range("P34").Value = [double]"2.6015445989699999"
$ws.Range("Q34").Value = [double]"1.913809067194"
$ws.Range("R34").Value = [double]"101.699215408141"
$ws.Range("S34").Value = [double]"3.1573596329470002"
$ws.Range("T34").Value = [double]"5.0674222287559996"
$ws.Range("U34").Value = [double]"2.1284911164799998"
$ws.Range("V34").Value = [double]"0.384692957329"
$ws.Range("W34").Value = [double]"0.66219819775599997"
$ws.Range("X34").Value = [double]"0.86448679854099997"

# Row 35
$ws.Range("N35").Value = [double]"3.5637682482999998E-2"
$ws.Range("O35").Value = [double]"3.3386797040999998E-2"
$ws.Range("P35").Value = [double]"2.3266583966060002"
$ws.Range("Q35").Value = [double]"1.810885146607"
$ws.Range("R35").Value = [double]"12.058820928453001"
$ws.Range("S35").Value = [double]"0.80665011879799997"
$ws.Range("T35").Value = [double]"1.203646389075"
$ws.Range("U35").Value = [double]"3.7320225332870001"
$ws.Range("V35").Value = [double]"0.30381714903599999"
$ws.Range("W35").Value = [double]"0.47990762164700002"
$ws.Range("X35").Value = [double]"0.62255663553899998"

# Row 36
$ws.Range("N36").Value = [double]"4.7972467957999997E-2"
$ws.Range("O36").Value = [double]"4.4850887919999999E-2"
$ws.Range("P36").Value = [double]"2.250652263389"
$ws.Range("Q36").Value = [double]"1.7857123390230001"
$ws.Range("R36").Value = [double]"4.4850608847059998"
$ws.Range("S36").Value = [double]"0.44566993641399999"
$ws.Range("T36").Value = [double]"0.65496093749999995"
$ws.Range("U36").Value = [double]"5.1240707945050001"
$ws.Range("V36").Value = [double]"0.27416447633500002"
$ws.Range("W36").Value = [double]"0.38830815155999998"
$ws.Range("X36").Value = [double]"0.50172450275199998"

# Row 37
$ws.Range("N37").Value = [double]"4.9328477629999997E-3"
$ws.Range("O37").Value = [double]"4.7701002950000003E-3"
$ws.Range("P37").Value = [double]"2.7200482714629999"
$ws.Range("Q37").Value = [double]"1.941020286973"
$ws.Range("R37").Value = [double]"513.86973280632299"
$ws.Range("S37").Value = [double]"14.808637434623"
$ws.Range("T37").Value = [double]"24.358715054867002"
$ws.Range("U37").Value = [double]"1.8386034632970001"
$ws.Range("V37").Value = [double]"0.43464075135000002"
$ws.Range("W37").Value = [double]"0.67814214752699997"
$ws.Range("X37").Value = [double]"0.88589416349200001"

# Row 38
$ws.Range("N38").Value = [double]"3.2197502854999997E-2"
$ws.Range("O38").Value = [double]"3.0922383768000002E-2"
$ws.Range("P38").Value = [double]"2.3957885727920001"
$ws.Range("Q38").Value = [double]"1.8341235743199999"
$ws.Range("R38").Value = [double]"36.514120132404003"
$ws.Range("S38").Value = [double]"2.0923396021910001"
$ws.Range("T38").Value = [double]"3.1393377807039999"
$ws.Range("U38").Value = [double]"3.0223645800660002"
$ws.Range("V38").Value = [double]"0.34546851162699999"
$ws.Range("W38").Value = [double]"0.52065539873699995"
$ws.Range("X38").Value = [double]"0.67628598570099996"

# Row 39
$ws.Range("N39").Value = [double]"4.5753599609000002E-2"
$ws.Range("O39").Value = [double]"4.3814211659999998E-2"
$ws.Range("P39").Value = [double]"2.2989321404590002"
$ws.Range("Q39").Value = [double]"1.801042356807"
$ws.Range("R39").Value = [double]"10.983961607873001"
$ws.Range("S39").Value = [double]"0.90216369038199995"
$ws.Range("T39").Value = [double]"1.3170759723600001"
$ws.Range("U39").Value = [double]"4.0371489632369997"
$ws.Range("V39").Value = [double]"0.30870345892899997"
$ws.Range("W39").Value = [double]"0.437417046985"
$ws.Range("X39").Value = [double]"0.56548651864800004"

# Row 40
$ws.Range("N40").Value = [double]"2.8175569740000002E-3"
$ws.Range("O40").Value = [double]"2.8266037849999999E-3"
$ws.Range("P40").Value = [double]"2.8735446417819999"
$ws.Range("Q40").Value = [double]"1.9722087181700001"
$ws.Range("R40").Value = [double]"4916.6786699054301"
$ws.Range("S40").Value = [double]"132.92945143268801"
$ws.Range("T40").Value = [double]"225.40543365478501"
$ws.Range("U40").Value = [double]"1.6247578148399999"
$ws.Range("V40").Value = [double]"0.48408474680199998"
$ws.Range("W40").Value = [double]"0.68866379871100003"
$ws.Range("X40").Value = [double]"0.90014898161500001"

# Row 41
$ws.Range("N41").Value = [double]"2.7243058779E-2"
$ws.Range("O41").Value = [double]"2.7102555702999999E-2"
$ws.Range("P41").Value = [double]"2.490143377011"
$ws.Range("Q41").Value = [double]"1.8649202338069999"
$ws.Range("R41").Value = [double]"162.42864808014099"
$ws.Range("S41").Value = [double]"8.0666463323070001"
$ws.Range("T41").Value = [double]"12.309307266375001"
$ws.Range("U41").Value = [double]"2.4936601678480002"
$ws.Range("V41").Value = [double]"0.39186810973500003"
$ws.Range("W41").Value = [double]"0.55569105681800002"
$ws.Range("X41").Value = [double]"0.72315880768100005"

# Row 42
$ws.Range("N42").Value = [double]"4.1814717132000002E-2"
$ws.Range("O42").Value = [double]"4.1413592397000001E-2"
$ws.Range("P42").Value = [double]"2.3772578897379999"
$ws.Range("Q42").Value = [double]"1.831500762611"
$ws.Range("R42").Value = [double]"35.753100829746998"
$ws.Range("S42").Value = [double]"2.4733025538760001"
$ws.Range("T42").Value = [double]"3.6590898437499999"
$ws.Range("U42").Value = [double]"3.2187576969810001"
$ws.Range("V42").Value = [double]"0.351697841751"
$ws.Range("W42").Value = [double]"0.48087345153599997"
$ws.Range("X42").Value = [double]"0.623181295977"

# Row 43
$ws.Range("N43").Value = [double]"5.3859023610000004E-3"
$ws.Range("O43").Value = [double]"5.076642349E-3"
$ws.Range("P43").Value = [double]"2.8165990491790001"
$ws.Range("Q43").Value = [double]"1.9742936897159999"
$ws.Range("R43").Value = [double]"89.950986081406995"
$ws.Range("S43").Value = [double]"2.557027428269"
$ws.Range("T43").Value = [double]"4.6624886718749998"
$ws.Range("U43").Value = [double]"2.1462201203419999"
$ws.Range("V43").Value = [double]"0.36845872756600001"
$ws.Range("W43").Value = [double]"0.68665685994299996"
$ws.Range("X43").Value = [double]"0.89950988068000004"

# Row 44
$ws.Range("N44").Value = [double]"3.0396316626000001E-2"
$ws.Range("O44").Value = [double]"2.8509465267999998E-2"
$ws.Range("P44").Value = [double]"2.5353751963039999"
$ws.Range("Q44").Value = [double]"1.9053130278929999"
$ws.Range("R44").Value = [double]"11.379634341960999"
$ws.Range("S44").Value = [double]"0.67122405056699996"
$ws.Range("T44").Value = [double]"1.094223983911"
$ws.Range("U44").Value = [double]"3.6867232972220001"
$ws.Range("V44").Value = [double]"0.29007437544999998"
$ws.Range("W44").Value = [double]"0.50783092639299998"
$ws.Range("X44").Value = [double]"0.66387242742700003"

# Row 45
$ws.Range("N45").Value = [double]"4.1884958894999998E-2"
$ws.Range("O45").Value = [double]"3.9212531597E-2"
$ws.Range("P45").Value = [double]"2.4419676921089999"
$ws.Range("Q45").Value = [double]"1.8811493437569999"
$ws.Range("R45").Value = [double]"4.3129056728539998"
$ws.Range("S45").Value = [double]"0.36912582408799999"
$ws.Range("T45").Value = [double]"0.58271982421900004"
$ws.Range("U45").Value = [double]"5.0310028176569999"
$ws.Range("V45").Value = [double]"0.25808238063400002"
$ws.Range("W45").Value = [double]"0.41883051948900002"
$ws.Range("X45").Value = [double]"0.54592907123800005"

# Row 46
$ws.Range("N46").Value = [double]"3.508901238E-3"
$ws.Range("O46").Value = [double]"3.3943331079999998E-3"
$ws.Range("P46").Value = [double]"2.9476790312559999"
$ws.Range("Q46").Value = [double]"1.986528250244"
$ws.Range("R46").Value = [double]"441.90020693285999"
$ws.Range("S46").Value = [double]"11.858146345844"
$ws.Range("T46").Value = [double]"22.623111246766001"
$ws.Range("U46").Value = [double]"1.858905884049"
$ws.Range("V46").Value = [double]"0.41826122403400001"
$ws.Range("W46").Value = [double]"0.697986161833"
$ws.Range("X46").Value = [double]"0.91485060835900001"

# Row 47
$ws.Range("N47").Value = [double]"2.7319876855000001E-2"
$ws.Range("O47").Value = [double]"2.6272990821999999E-2"
$ws.Range("P47").Value = [double]"2.6182657818599999"
$ws.Range("Q47").Value = [double]"1.9263368002229999"
$ws.Range("R47").Value = [double]"34.028458764753999"
$ws.Range("S47").Value = [double]"1.75424347529"
$ws.Range("T47").Value = [double]"2.9141839443549999"
$ws.Range("U47").Value = [double]"2.9902422627190002"
$ws.Range("V47").Value = [double]"0.33364927362800001"
$ws.Range("W47").Value = [double]"0.54390210430099994"
$ws.Range("X47").Value = [double]"0.71189369453999995"

# Row 48
$ws.Range("N48").Value = [double]"4.0071067342999998E-2"
$ws.Range("O48").Value = [double]"3.8433482398999998E-2"
$ws.Range("P48").Value = [double]"2.5127986323860001"
$ws.Range("Q48").Value = [double]"1.901276080785"
$ws.Range("R48").Value = [double]"10.466873692492999"
$ws.Range("S48").Value = [double]"0.76016330623700001"
$ws.Range("T48").Value = [double]"1.210736980469"
$ws.Range("U48").Value = [double]"3.9587947162030002"
$ws.Range("V48").Value = [double]"0.29646678164200002"
$ws.Range("W48").Value = [double]"0.462968039981"
$ws.Range("X48").Value = [double]"0.60412399437999997"

# Row 49
$ws.Range("N49").Value = [double]"1.916153124E-3"
$ws.Range("O49").Value = [double]"1.922789668E-3"
$ws.Range("P49").Value = [double]"3.104447332791"
$ws.Range("Q49").Value = [double]"1.990770069114"
$ws.Range("R49").Value = [double]"4063.9004758534102"
$ws.Range("S49").Value = [double]"103.605693752908"
$ws.Range("T49").Value = [double]"209.54890828363301"
$ws.Range("U49").Value = [double]"1.648085053877"
$ws.Range("V49").Value = [double]"0.46668341023400001"
$ws.Range("W49").Value = [double]"0.70518159165299998"
$ws.Range("X49").Value = [double]"0.92470402135999996"

# Row 50
$ws.Range("N50").Value = [double]"2.2616050897000001E-2"
$ws.Range("O50").Value = [double]"2.2538056808000001E-2"
$ws.Range("P50").Value = [double]"2.7240910689070001"
$ws.Range("Q50").Value = [double]"1.9485570637159999"
$ws.Range("R50").Value = [double]"148.95090733415"
$ws.Range("S50").Value = [double]"6.7070100998570004"
$ws.Range("T50").Value = [double]"11.37802381"
$ws.Range("U50").Value = [double]"2.4828945102510001"
$ws.Range("V50").Value = [double]"0.37833181377000002"
$ws.Range("W50").Value = [double]"0.57721824766200003"
$ws.Range("X50").Value = [double]"0.75623262736600005"

# Row 51
$ws.Range("N51").Value = [double]"3.5857065969999997E-2"
$ws.Range("O51").Value = [double]"3.5593186994999998E-2"
$ws.Range("P51").Value = [double]"2.5966152982409998"
$ws.Range("Q51").Value = [double]"1.923989066828"
$ws.Range("R51").Value = [double]"33.689782715881002"
$ws.Range("S51").Value = [double]"2.057957182425"
$ws.Range("T51").Value = [double]"3.298362890625"
$ws.Range("U51").Value = [double]"3.186816635449"
$ws.Range("V51").Value = [double]"0.33577339530900002"
$ws.Range("W51").Value = [double]"0.507393625811"
$ws.Range("X51").Value = [double]"0.662550642488"

# Restore the view: scroll the frozen pane back to show column B at the top
# (previously scrolled right to column G), and clear the explicit selection
# in the bottom-right pane.
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 2
